# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to match newly scraped counts (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 970
$ws1.Range("F6").Value  = 5304
$ws1.Range("F8").Value  = 639
$ws1.Range("F9").Value  = 918
$ws1.Range("F10").Value = 819
$ws1.Range("F14").Value = 18
$ws1.Range("F17").Value = 1773
$ws1.Range("F18").Value = 1455
$ws1.Range("F19").Value = 838
$ws1.Range("F22").Value = 310
$ws1.Range("F23").Value = 518
$ws1.Range("F25").Value = 1047
$ws1.Range("F28").Value = 2624
$ws1.Range("F34").Value = 280
$ws1.Range("F39").Value = 275

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 970
$ws4.Range("F7").Value  = 5304
$ws4.Range("F9").Value  = 639
$ws4.Range("F12").Value = 918
$ws4.Range("F13").Value = 819
$ws4.Range("F19").Value = 18
$ws4.Range("F23").Value = 1773
$ws4.Range("F24").Value = 1455
$ws4.Range("F25").Value = 838
$ws4.Range("F27").Value = 310
$ws4.Range("F29").Value = 518
$ws4.Range("F31").Value = 1047
$ws4.Range("F33").Value = 2624
$ws4.Range("F38").Value = 280
$ws4.Range("F42").Value = 275
